$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "Total Buy Signals" header in G1, matching the bold/
# bordered/centered style already used by the other header cells (copy
# format from F1, then overwrite the value/text). ---
$ws.Range("F1").Copy($ws.Range("G1"))
$ws.Cells.Item(1, 7).Value2 = "Total Buy Signals"

# --- Rewrite the data rows (2-11) with the refreshed backtest numbers and
# the new "Total Buy Signals" column (G). The A (Stock Ticker), B
# (Algorithm) and C (methodology) text values are re-sent even where they
# are unchanged, and D (Initial Balance) is re-sent too, just to keep each
# row self-contained and correct. ---

# Row 2: NFE / Ensemble / classification
$ws.Cells.Item(2, 1).Value2 = "NFE"
$ws.Cells.Item(2, 2).Value2 = "Ensemble"
$ws.Cells.Item(2, 3).Value2 = "classification"
$ws.Cells.Item(2, 4).Value2 = 1000
$ws.Cells.Item(2, 5).Value2 = -168.3878771182665
$ws.Cells.Item(2, 6).Value2 = -116.8387877118267
$ws.Cells.Item(2, 7).Value2 = 249

# Row 3: NFE / LSTM / classification
$ws.Cells.Item(3, 1).Value2 = "NFE"
$ws.Cells.Item(3, 2).Value2 = "LSTM"
$ws.Cells.Item(3, 3).Value2 = "classification"
$ws.Cells.Item(3, 4).Value2 = 1000
$ws.Cells.Item(3, 5).Value2 = 239.1518985677852
$ws.Cells.Item(3, 6).Value2 = -76.08481014322149
$ws.Cells.Item(3, 7).Value2 = 203

# Row 4: NFE / Ensemble / regression
$ws.Cells.Item(4, 1).Value2 = "NFE"
$ws.Cells.Item(4, 2).Value2 = "Ensemble"
$ws.Cells.Item(4, 3).Value2 = "regression"
$ws.Cells.Item(4, 4).Value2 = 1000
$ws.Cells.Item(4, 5).Value2 = 1000
$ws.Cells.Item(4, 6).Value2 = 0
$ws.Cells.Item(4, 7).Value2 = 0

# Row 5: NFE / LSTM / regression
$ws.Cells.Item(5, 1).Value2 = "NFE"
$ws.Cells.Item(5, 2).Value2 = "LSTM"
$ws.Cells.Item(5, 3).Value2 = "regression"
$ws.Cells.Item(5, 4).Value2 = 1000
$ws.Cells.Item(5, 5).Value2 = 467.5304207126891
$ws.Cells.Item(5, 6).Value2 = -53.24695792873109
$ws.Cells.Item(5, 7).Value2 = 97

# Row 6: NFE / VWAP (no methodology)
$ws.Cells.Item(6, 1).Value2 = "NFE"
$ws.Cells.Item(6, 2).Value2 = "VWAP"
$ws.Cells.Item(6, 4).Value2 = 1000
$ws.Cells.Item(6, 5).Value2 = 610.707381687405
$ws.Cells.Item(6, 6).Value2 = -38.9292618312595
$ws.Cells.Item(6, 7).Value2 = 13

# Row 7: PLUG / Ensemble / classification
$ws.Cells.Item(7, 1).Value2 = "PLUG"
$ws.Cells.Item(7, 2).Value2 = "Ensemble"
$ws.Cells.Item(7, 3).Value2 = "classification"
$ws.Cells.Item(7, 4).Value2 = 1000
$ws.Cells.Item(7, 5).Value2 = 697.2485238905774
$ws.Cells.Item(7, 6).Value2 = -30.27514761094226
$ws.Cells.Item(7, 7).Value2 = 249

# Row 8: PLUG / LSTM / classification
$ws.Cells.Item(8, 1).Value2 = "PLUG"
$ws.Cells.Item(8, 2).Value2 = "LSTM"
$ws.Cells.Item(8, 3).Value2 = "classification"
$ws.Cells.Item(8, 4).Value2 = 1000
$ws.Cells.Item(8, 5).Value2 = 697.2485238905774
$ws.Cells.Item(8, 6).Value2 = -30.27514761094226
$ws.Cells.Item(8, 7).Value2 = 249

# Row 9: PLUG / Ensemble / regression
$ws.Cells.Item(9, 1).Value2 = "PLUG"
$ws.Cells.Item(9, 2).Value2 = "Ensemble"
$ws.Cells.Item(9, 3).Value2 = "regression"
$ws.Cells.Item(9, 4).Value2 = 1000
$ws.Cells.Item(9, 5).Value2 = 1000
$ws.Cells.Item(9, 6).Value2 = 0
$ws.Cells.Item(9, 7).Value2 = 0

# Row 10: PLUG / LSTM / regression
$ws.Cells.Item(10, 1).Value2 = "PLUG"
$ws.Cells.Item(10, 2).Value2 = "LSTM"
$ws.Cells.Item(10, 3).Value2 = "regression"
$ws.Cells.Item(10, 4).Value2 = 1000
$ws.Cells.Item(10, 5).Value2 = 1313.491649519817
$ws.Cells.Item(10, 6).Value2 = 31.34916495198169
$ws.Cells.Item(10, 7).Value2 = 31

# Row 11: PLUG / VWAP (no methodology)
$ws.Cells.Item(11, 1).Value2 = "PLUG"
$ws.Cells.Item(11, 2).Value2 = "VWAP"
$ws.Cells.Item(11, 4).Value2 = 1000
$ws.Cells.Item(11, 5).Value2 = 389.1936620642327
$ws.Cells.Item(11, 6).Value2 = -61.08063379357672
$ws.Cells.Item(11, 7).Value2 = 21
